$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a fresh row below the current last data row (19) to hold the new period.
$ws.Rows("20").Insert()

# 2) Give the new row 20 the same "closing" look (bottom border, etc.) that row 19
#    currently has, along with its values, by copying row 19 into row 20.
$ws.Range("B19:J19").Copy($ws.Range("B20:J20"))

# 3) Row 19 is no longer the last row of the table, so restyle it like the interior
#    rows above it (16-18) by copying formatting + values from row 18, then fix the
#    period back to 2508.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$ws.Range("E19").Value = "2508"

# 4) The newly duplicated row 20 still has the old "2508" period - bump it to 2509.
$ws.Range("E20").Value = "2509"

# 5) Update the summary figures to reflect the extra period.
$ws.Range("E11").Value = 244842
$ws.Range("F13").Value = 5

$excel.CutCopyMode = $false
